$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 252
$ws.Range("I2").Value = 240
$ws.Range("K2").Value = 240
$ws.Range("M2").Value = -127

$ws.Range("H17").Value = 2500891.5
$ws.Range("J17").Value = 2500891.5
$ws.Range("L17").Value = 7502674.5
$ws.Range("N17").Value = -7503010.5

$ws.Range("H62").Value = 3305.3704
$ws.Range("I62").Value = 2759.3157
$ws.Range("J62").Value = 4602.25
$ws.Range("K62").Value = 2759.3157
$ws.Range("L62").Value = 4602.25
$ws.Range("M62").Value = -2135.3157
$ws.Range("N62").Value = -5850.25

$ws.Range("H65").Value = 3305.3704
$ws.Range("I65").Value = 2759.3157
$ws.Range("J65").Value = 4602.25
$ws.Range("K65").Value = 13796.5785
$ws.Range("L65").Value = 23011.25
$ws.Range("M65").Value = -10676.5785
$ws.Range("N65").Value = -29251.25

$ws.Range("H112").Value = 1068.36
$ws.Range("J112").Value = 1096.2609
$ws.Range("L112").Value = 3288.7827
$ws.Range("N112").Value = -5504.7827

$ws.Range("H116").Value = 3753.7144
$ws.Range("I116").Value = 4584.2856
$ws.Range("J116").Value = 2923.1428
$ws.Range("K116").Value = 4584.2856
$ws.Range("L116").Value = 2923.1428
$ws.Range("M116").Value = -1142.2856
$ws.Range("N116").Value = -9807.1428

$ws.Range("H138").Value = 2141.2354
$ws.Range("I138").Value = 1579.921
$ws.Range("K138").Value = 4739.763
$ws.Range("M138").Value = 400.2370000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1467.76
$ws.Range("I45").Value = 1434
$ws.Range("J45").Value = 1602.8
$ws.Range("K45").Value = 1434
$ws.Range("L45").Value = 1602.8
$ws.Range("M45").Value = -1057
$ws.Range("N45").Value = -2356.8

$ws.Range("H132").Value = 2013.6471
$ws.Range("I132").Value = 1345.52
$ws.Range("K132").Value = 4036.56
$ws.Range("M132").Value = -1506.56

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1604.4375
$ws.Range("I20").Value = 1435.1578
$ws.Range("J20").Value = 1851.8462
$ws.Range("K20").Value = 1435.1578
$ws.Range("L20").Value = 1851.8462
$ws.Range("M20").Value = -1188.1578
$ws.Range("N20").Value = -2345.8462

$ws.Range("H133").Value = 39788.445
$ws.Range("I133").Value = 30709
$ws.Range("J133").Value = 40923.375
$ws.Range("K133").Value = 30709
$ws.Range("L133").Value = 40923.375
$ws.Range("M133").Value = -25649
$ws.Range("N133").Value = -51043.375

$ws.Range("H140").Value = 42512.223
$ws.Range("J140").Value = 42512.223
$ws.Range("L140").Value = 42512.223
$ws.Range("N140").Value = -52872.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 973.8461
$ws.Range("I16").Value = 895
$ws.Range("J16").Value = 1236.6666
$ws.Range("K16").Value = 895
$ws.Range("L16").Value = 1236.6666
$ws.Range("M16").Value = -608
$ws.Range("N16").Value = -1810.6666

$ws.Range("H99").Value = 2919971.5
$ws.Range("J99").Value = 27000
$ws.Range("L99").Value = 27000
$ws.Range("N99").Value = -29996

$ws.Range("H113").Value = 973.8461
$ws.Range("I113").Value = 895
$ws.Range("J113").Value = 1236.6666
$ws.Range("K113").Value = 895
$ws.Range("L113").Value = 1236.6666
$ws.Range("M113").Value = 1275
$ws.Range("N113").Value = -5576.6666

$ws.Range("H126").Value = 2919971.5
$ws.Range("J126").Value = 27000
$ws.Range("L126").Value = 81000
$ws.Range("N126").Value = -85940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3353.4546
$ws.Range("I94").Value = 595
$ws.Range("J94").Value = 3966.4443
$ws.Range("K94").Value = 1785
$ws.Range("L94").Value = 11899.3329
$ws.Range("M94").Value = -1109
$ws.Range("N94").Value = -13251.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5440.394
$ws.Range("I70").Value = 4973.95
$ws.Range("K70").Value = 4973.95
$ws.Range("M70").Value = -4703.95

$ws.Range("H73").Value = 5440.394
$ws.Range("I73").Value = 4973.95
$ws.Range("K73").Value = 4973.95
$ws.Range("M73").Value = -4037.95

$ws.Range("H97").Value = 612.0714
$ws.Range("I97").Value = 595.9
$ws.Range("J97").Value = 652.5
$ws.Range("K97").Value = 595.9
$ws.Range("L97").Value = 652.5
$ws.Range("M97").Value = -99.89999999999998
$ws.Range("N97").Value = -1644.5

$ws.Range("H102").Value = 3298.1875
$ws.Range("I102").Value = 2975.7
$ws.Range("J102").Value = 3835.6667
$ws.Range("K102").Value = 2975.7
$ws.Range("L102").Value = 3835.6667
$ws.Range("M102").Value = -1353.7
$ws.Range("N102").Value = -7079.6667

$ws.Range("H113").Value = 1656.3572
$ws.Range("I113").Value = 1518.9
$ws.Range("K113").Value = 1518.9
$ws.Range("M113").Value = 651.0999999999999

$ws.Range("H122").Value = 2179.3333
$ws.Range("I122").Value = 1655.1177
$ws.Range("K122").Value = 4965.3531
$ws.Range("M122").Value = -2515.3531

$ws.Range("H126").Value = 2266.2964
$ws.Range("I126").Value = 1922.1428
$ws.Range("J126").Value = 2636.923
$ws.Range("K126").Value = 5766.428400000001
$ws.Range("L126").Value = 7910.768999999999
$ws.Range("M126").Value = -3296.428400000001
$ws.Range("N126").Value = -12850.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 38307.75
$ws.Range("I7").Value = 54816.21
$ws.Range("J7").Value = 3456.5557
$ws.Range("K7").Value = 54816.21
$ws.Range("L7").Value = 3456.5557
$ws.Range("M7").Value = -54704.21
$ws.Range("N7").Value = -3680.5557

$ws.Range("H93").Value = 7725
$ws.Range("I93").Value = 9631.454
$ws.Range("J93").Value = 734.6667
$ws.Range("K93").Value = 9631.454
$ws.Range("L93").Value = 734.6667
$ws.Range("M93").Value = -8383.454
$ws.Range("N93").Value = -3230.6667

$ws.Range("H100").Value = 1109.8572
$ws.Range("I100").Value = 1109.8572
$ws.Range("K100").Value = 1109.8572
$ws.Range("M100").Value = -568.8571999999999

$ws.Range("H122").Value = 2531.4285
$ws.Range("I122").Value = 1918.5
$ws.Range("J122").Value = 2991.125
$ws.Range("K122").Value = 5755.5
$ws.Range("L122").Value = 8973.375
$ws.Range("M122").Value = -3305.5
$ws.Range("N122").Value = -13873.375

$ws.Range("H126").Value = 38307.75
$ws.Range("I126").Value = 54816.21
$ws.Range("J126").Value = 3456.5557
$ws.Range("K126").Value = 164448.63
$ws.Range("L126").Value = 10369.6671
$ws.Range("M126").Value = -161978.63
$ws.Range("N126").Value = -15309.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 56954.5
$ws.Range("I122").Value = 72484.36
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 217453.08
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -215003.08
$ws.Range("N122").Value = -12700

$ws.Range("H126").Value = 126137.25
$ws.Range("I126").Value = 143599
$ws.Range("J126").Value = 3905
$ws.Range("K126").Value = 430797
$ws.Range("L126").Value = 11715
$ws.Range("M126").Value = -428327
$ws.Range("N126").Value = -16655

$ws.Range("H139").Value = 70427.78
$ws.Range("J139").Value = 78571.42999999999
$ws.Range("L139").Value = 78571.42999999999
$ws.Range("N139").Value = -88851.42999999999
